# Update Sheets via scheduled runner
# Applies refreshed market data values to the Anima Profits workbook
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 1026.8572
$ws.Range("I34").Value = 1026.8572
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1026.8572
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -823.8571999999999
$ws.Range("N34").ClearContents()

# Row 36
$ws.Range("H36").Value = 1026.8572
$ws.Range("I36").Value = 1026.8572
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1026.8572
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -311.8571999999999
$ws.Range("N36").ClearContents()

# Row 106
$ws.Range("H106").Value = 6668635
$ws.Range("I106").Value = 7144680.5
$ws.Range("K106").Value = 7144680.5
$ws.Range("M106").Value = -7144049.5
$ws.Range("N106").ClearContents()

# Row 113
$ws.Range("H113").Value = 2838.3845
$ws.Range("I113").Value = 2700
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = 554
$ws.Range("N113").ClearContents()

# Row 132
$ws.Range("H132").Value = 2269.775
$ws.Range("I132").Value = 2122.9473
$ws.Range("J132").Value = 5059.5
$ws.Range("K132").Value = 6368.841899999999
$ws.Range("L132").Value = 15178.5
$ws.Range("M132").Value = -3838.841899999999
$ws.Range("N132").Value = -20238.5

# Row 133
$ws.Range("H133").Value = 57585
$ws.Range("J133").Value = 57585
$ws.Range("L133").Value = 57585
$ws.Range("N133").Value = -67705

# Row 136
$ws.Range("H136").Value = 32585
$ws.Range("J136").Value = 32585
$ws.Range("L136").Value = 32585
$ws.Range("N136").Value = -42785

# Row 138
$ws.Range("H138").Value = 2217.3684
$ws.Range("I138").Value = 1839.9615
$ws.Range("J138").Value = 3035.0833
$ws.Range("K138").Value = 5519.8845
$ws.Range("L138").Value = 9105.249899999999
$ws.Range("M138").Value = -379.8845000000001
$ws.Range("N138").Value = -19385.2499

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 564453.7
$ws.Range("I32").Value = 678534.3
$ws.Range("K32").Value = 678534.3
$ws.Range("M32").Value = -678247.3
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 27
$ws.Range("H27").Value = 69990
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

# Row 123
$ws.Range("H123").Value = 50675
$ws.Range("J123").Value = 50675
$ws.Range("L123").Value = 50675
$ws.Range("N123").Value = -60475

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 46614.832
$ws.Range("I23").Value = 3834.5
$ws.Range("J23").Value = 68005
$ws.Range("K23").Value = 3834.5
$ws.Range("L23").Value = 68005
$ws.Range("M23").Value = -3594.5
$ws.Range("N23").Value = -68485

# Row 27
$ws.Range("H27").Value = 46614.832
$ws.Range("I27").Value = 3834.5
$ws.Range("J27").Value = 68005
$ws.Range("K27").Value = 3834.5
$ws.Range("L27").Value = 68005
$ws.Range("M27").Value = -3642.5
$ws.Range("N27").Value = -68389

# Row 31
$ws.Range("H31").Value = 5862.32
$ws.Range("I31").Value = 1305.8928
$ws.Range("J31").Value = 11661.409
$ws.Range("K31").Value = 1305.8928
$ws.Range("L31").Value = 11661.409
$ws.Range("M31").Value = -1010.8928
$ws.Range("N31").Value = -12251.409

# Row 34
$ws.Range("H34").Value = 5862.32
$ws.Range("I34").Value = 1305.8928
$ws.Range("J34").Value = 11661.409
$ws.Range("K34").Value = 1305.8928
$ws.Range("L34").Value = 11661.409
$ws.Range("M34").Value = -1103.8928
$ws.Range("N34").Value = -12065.409

# Row 58
$ws.Range("H58").Value = 1994.95
$ws.Range("I58").Value = 1813.6
$ws.Range("K58").Value = 1813.6
$ws.Range("M58").Value = -1610.6
$ws.Range("N58").ClearContents()

# Row 136
$ws.Range("H136").Value = 1994.95
$ws.Range("I136").Value = 1813.6
$ws.Range("K136").Value = 5440.799999999999
$ws.Range("M136").Value = -2890.799999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 9196
$ws.Range("J33").Value = 100
$ws.Range("L33").Value = 600
$ws.Range("N33").Value = -1166

# Row 40
$ws.Range("H40").Value = 137.6
$ws.Range("I40").Value = 117.70588
$ws.Range("J40").Value = 179.875
$ws.Range("K40").Value = 470.82352
$ws.Range("L40").Value = 719.5
$ws.Range("M40").Value = -401.82352
$ws.Range("N40").Value = -857.5

# Row 118
$ws.Range("H118").Value = 2966.1538
$ws.Range("J118").Value = 3017.1428
$ws.Range("L118").Value = 9051.428400000001
$ws.Range("N118").Value = -11537.4284

# Row 122
$ws.Range("H122").Value = 5075.826
$ws.Range("I122").Value = 403.07144
$ws.Range("J122").Value = 12344.556
$ws.Range("K122").Value = 3627.64296
$ws.Range("L122").Value = 111101.004
$ws.Range("M122").Value = -1177.64296
$ws.Range("N122").Value = -116001.004

# Row 126
$ws.Range("H126").Value = 5146.2
$ws.Range("I126").Value = 1799
$ws.Range("J126").Value = 5385.2856
$ws.Range("K126").Value = 5397
$ws.Range("L126").Value = 16155.8568
$ws.Range("M126").Value = -457
$ws.Range("N126").Value = -26035.8568

# Row 131
$ws.Range("H131").Value = 1023.55817
$ws.Range("I131").Value = 631.2857
$ws.Range("J131").Value = 1099.8334
$ws.Range("K131").Value = 1893.8571
$ws.Range("L131").Value = 3299.5002
$ws.Range("M131").Value = 3146.1429
$ws.Range("N131").Value = -13379.5002

# Row 132
$ws.Range("H132").Value = 2264.3635
$ws.Range("I132").Value = 2195.7273
$ws.Range("J132").Value = 2287.2424
$ws.Range("K132").Value = 19761.5457
$ws.Range("L132").Value = 20585.1816
$ws.Range("M132").Value = -17231.5457
$ws.Range("N132").Value = -25645.1816

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 6799735
$ws.Range("I10").Value = 15012500
$ws.Range("J10").Value = 229523.2
$ws.Range("K10").Value = 15012500
$ws.Range("L10").Value = 229523.2
$ws.Range("M10").Value = -15012331
$ws.Range("N10").Value = -229861.2

# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 102
$ws.Range("H102").Value = 2400
$ws.Range("I102").Value = 2350
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2350
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -728
$ws.Range("N102").Value = -5744

# Row 122
$ws.Range("H122").Value = 2429.7144
$ws.Range("I122").Value = 1669.3334
$ws.Range("K122").Value = 5008.0002
$ws.Range("M122").Value = -2558.0002
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 3179.5588
$ws.Range("I132").Value = 2684.4
$ws.Range("J132").Value = 4555
$ws.Range("K132").Value = 8053.200000000001
$ws.Range("L132").Value = 13665
$ws.Range("M132").Value = -5523.200000000001
$ws.Range("N132").Value = -18725

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2233.9644
$ws.Range("I132").Value = 1575.7
$ws.Range("J132").Value = 3879.625
$ws.Range("K132").Value = 4727.1
$ws.Range("L132").Value = 11638.875
$ws.Range("M132").Value = -2197.1
$ws.Range("N132").Value = -16698.875

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 46403
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 57503.75
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 57503.75
$ws.Range("M7").Value = -1887
$ws.Range("N7").Value = -57729.75

# Row 51
$ws.Range("H51").Value = 14250
$ws.Range("J51").Value = 14250
$ws.Range("L51").Value = 14250
$ws.Range("N51").Value = -15270

# Row 136
$ws.Range("H136").Value = 2062.4583
$ws.Range("I136").Value = 1505.3611
$ws.Range("K136").Value = 4516.0833
$ws.Range("M136").Value = -1966.0833
$ws.Range("N136").ClearContents()
